# Insert a new weekly price record at row 269 for "Pepino dulce" / Vega Modelo
# de Temuco. All existing rows from 269 down to 313 shift down by one (to
# 270-314), preserving their original data, and the sheet's used range grows
# by one row (A1:R313 -> A1:R314).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 269..313 down to 270..314, leaving a fresh blank row at 269.
$ws.Rows(269).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(269, 1).Value = 10
$ws.Cells.Item(269, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(269, 3).Value = "La Araucanía"
$ws.Cells.Item(269, 4).Value = 44995
$ws.Cells.Item(269, 5).Value = 9
$ws.Cells.Item(269, 6).Value = 100112043
$ws.Cells.Item(269, 7).Value = "Pepino dulce"
$ws.Cells.Item(269, 8).Value = "Cultivar XV región"
$ws.Cells.Item(269, 9).Value = "Primera"
$ws.Cells.Item(269, 10).Value = 35
$ws.Cells.Item(269, 11).Value = 17000
$ws.Cells.Item(269, 12).Value = 17000
$ws.Cells.Item(269, 13).Value = 17000
$ws.Cells.Item(269, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(269, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(269, 16).Value = 944
$ws.Cells.Item(269, 17).Value = 18
$ws.Cells.Item(269, 18).Value = "Hortaliza"
